$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for every data row (2-294)
# from 45205 (2023-10-06) to 45206 (2023-10-07).
$ws.Range("C2:C294").Value = 45206
